$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RegistrationNo value change
$ws.Range("B3").Value = "SLB-RLOS-131580"

# Row 6: TRN, RegistrationNo and SemesterCost change
$ws.Range("A6").Value = 110732120
$ws.Range("B6").Value = "SLB-RLOS-129449"
$ws.Range("D6").Value = 516608.67999999993

# Old row 8 is removed entirely (its content no longer appears anywhere in the sheet)
$ws.Range("A8:H8").ClearContents()

# Row 9 gets a brand new set of values (not simply shifted from the old row 10)
$ws.Range("A9").Value = 127645999
$ws.Range("B9").Value = "SLB-RLOS-141000"
$ws.Range("C9").Value = "YES"
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0
$ws.Range("H9").Value = 2023

# Old row 10 is removed entirely (its content no longer appears anywhere in the sheet)
$ws.Range("A10:H10").ClearContents()
